# Ratio Matrix Construction "Done"
# Applies the recomputed ratio-matrix values to FacilityInfo / CO2LocationInfo.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FacilityInfo")
$ws2 = $wb.Worksheets.Item("CO2LocationInfo")

# --- FacilityInfo: re-ordered usage columns AF:AJ (header + value move together) ---
# Headers (shared-string table reorder in the source diff; net effect is these
# five header cells now show the labels in this column order).
$ws1.Range("AF1").Value = "tkm-N1Usage"
$ws1.Range("AG1").Value = "pkmUsage"
$ws1.Range("AH1").Value = "tkm-SZMUsage"
$ws1.Range("AI1").Value = "tkm-N2Usage"
$ws1.Range("AJ1").Value = "tkm-N3Usage"

# Matching data row, reordered the same way as its header.
$ws1.Range("AF2").Value2 = 7.5
$ws1.Range("AG2").Value2 = 850
$ws1.Range("AH2").Value2 = 414.5
$ws1.Range("AI2").Value2 = 24.2
$ws1.Range("AJ2").Value2 = 130.3

# --- FacilityInfo: recomputed values on row 2 ---
$ws1.Range("B2").Value2 = 60374617778.56241
$ws1.Range("F2").Value2 = 343.4814507244961
$ws1.Range("N2").Value2 = 288.5244186085768
$ws1.Range("O2").Value2 = 288.5244186085768
$ws1.Range("P2").Value2 = -0.000000000000007553957459549564
$ws1.Range("Q2").Value2 = -0.000000000000007553957459549564

# --- CO2LocationInfo: recomputed "Amount Used" (column D) ---
$ws2.Range("D2").Value2 = 1.483469072164948
$ws2.Range("D4").Value2 = 6.790639175257732
$ws2.Range("D5").Value2 = 7.718798969072166
$ws2.Range("D17").Value2 = 1.769056701030928
$ws2.Range("D19").Value2 = 1.784922680412372
$ws2.Range("D22").Value2 = 3.879231958762886
$ws2.Range("D29").Value2 = 0.8726288659793813
$ws2.Range("D35").Value2 = 0.9598917525773195
$ws2.Range("D36").Value2 = 1.301010309278351
$ws2.Range("D45").Value2 = 0.9598917525773195
$ws2.Range("D50").Value2 = 6.39398969072165
$ws2.Range("D56").Value2 = 4.910520618556702
$ws2.Range("D60").Value2 = 5.545159793814434
$ws2.Range("D64").Value2 = 4.505938144329898
$ws2.Range("D66").Value2 = 7.528407216494846
$ws2.Range("D71").Value2 = 6.39398969072165
$ws2.Range("D80").Value2 = 1.388273195876288
$ws2.Range("D81").Value2 = 1.665927835051547
$ws2.Range("D83").Value2 = 1.967381443298969
$ws2.Range("D91").Value2 = 1.348608247422681
$ws2.Range("D94").Value2 = 6.917567010309281
$ws2.Range("D98").Value2 = 0.9519587628865978
$ws2.Range("D104").Value2 = 2.35609793814433
$ws2.Range("D105").Value2 = 1.761123711340206
$ws2.Range("D108").Value2 = 0.9678247422680412
$ws2.Range("D113").Value2 = 0.8726288659793813
$ws2.Range("D117").Value2 = 1.697659793814433
$ws2.Range("D124").Value2 = 0.9995567010309281
$ws2.Range("D131").Value2 = 1.753190721649485
$ws2.Range("D132").Value2 = 3.187924951300077
